$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("Contents")
Write-Output ("Start=" + $rng.Start + " End=" + $rng.End)
try {
  $rng.InsertAfter("ZZZ")
  Write-Output "InsertAfter ok"
} catch {
  Write-Output ("ERR InsertAfter: " + $_)
}
